# Actualización automática 2025-06-21 14:30:08
#
# Registers a new sale for HIDALGO HIDALGO PEDRO GUSTAVO / MEGAMAFERS S.A.
# in the PORCELANATO group for the current month (junio), and propagates
# the resulting totals through the dependent summary sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "VENTAS POR GRUPO" -------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M12").Value = 3228.67
$wsGrupo.Range("M22").Value = "3 de 20"

# --- Sheet 2: "VENTA MENSUAL" -----------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F12").Value = 3228.67
$wsMensual.Range("F22").Value = 3397.23

# --- Sheet 3: "CUMPLIMIENTO MENSUAL" ----------------------------------
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Range("D16").Value = 3451.99
$wsCumplimiento.Range("E16").Value = 26080.45
$wsCumplimiento.Range("F16").Value = 0.1168880729123635

$wsCumplimiento.Range("D19").Value = 3397.23
$wsCumplimiento.Range("E19").Value = 46989.96762291768
$wsCumplimiento.Range("F19").Value = 0.0674224834932045
